# Update the answer table: replace the computed division results in the
# data rows (1, 5, 9, 13, 17) of the single 5-column table with the new
# values from the regenerated output.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "27÷8=3, 3"
$t.Cell(1, 2).Range.Text  = "28÷5=5, 3"
$t.Cell(1, 3).Range.Text  = "14÷6=2, 2"
$t.Cell(1, 4).Range.Text  = "75÷4=18, 3"
$t.Cell(1, 5).Range.Text  = "97÷6=16, 1"

$t.Cell(5, 1).Range.Text  = "61÷2=30, 1"
$t.Cell(5, 2).Range.Text  = "75÷4=18, 3"
$t.Cell(5, 3).Range.Text  = "18÷9=2, 0"
$t.Cell(5, 4).Range.Text  = "41÷3=13, 2"
$t.Cell(5, 5).Range.Text  = "47÷6=7, 5"

$t.Cell(9, 1).Range.Text  = "17÷8=2, 1"
$t.Cell(9, 2).Range.Text  = "40÷2=20, 0"
$t.Cell(9, 3).Range.Text  = "71÷2=35, 1"
$t.Cell(9, 4).Range.Text  = "48÷4=12, 0"
$t.Cell(9, 5).Range.Text  = "25÷9=2, 7"

$t.Cell(13, 1).Range.Text = "69÷5=13, 4"
$t.Cell(13, 2).Range.Text = "13÷2=6, 1"
$t.Cell(13, 3).Range.Text = "24÷2=12, 0"
$t.Cell(13, 4).Range.Text = "16÷5=3, 1"
$t.Cell(13, 5).Range.Text = "51÷7=7, 2"

$t.Cell(17, 1).Range.Text = "13÷4=3, 1"
$t.Cell(17, 2).Range.Text = "48÷5=9, 3"
$t.Cell(17, 3).Range.Text = "81÷6=13, 3"
$t.Cell(17, 4).Range.Text = "29÷3=9, 2"
$t.Cell(17, 5).Range.Text = "85÷9=9, 4"
